$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("expected")

# "fixed string arguments in filters (breaking change)":
# the jinja-style filter template texts used an unquoted `Other` column
# name as an argument to maxrows(); quote it so it is treated as a
# string literal instead of a (non-existent) variable.
$ws1.Range("A3").Value = '{{ df2 | maxrows(2, "Other", 0)}}'
$ws1.Range("A6").Value = '{{ df2 | noheader | maxrows(2, "Other", 0) }}'

# Reflect the saved workbook's UI state: Sheet1 active with A7 selected.
$ws1.Activate()
$ws1.Range("A7").Select()
